# Fixed #418 Empty AQL expressions generate empty lines.
# Remove the empty paragraph (2nd paragraph) that was left behind, which
# contains only an empty run (w:t/>) and a tab stop paragraph property.

$d = $word.ActiveDocument

$d.Paragraphs(2).Range.Delete()
